# Weekly update: insert a new observation row for the Berenjena
# (eggplant) price sheet. The new record is inserted at row 104,
# pushing the existing rows 104-152 down to 105-153.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 104 - this shifts rows
# 104..152 down to 105..153 (and extends the used range to R153).
$ws.Rows.Item(104).Insert()

# Populate the newly inserted row 104 with the new weekly record.
$ws.Cells.Item(104, 1).Value = 6
$ws.Cells.Item(104, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(104, 3).Value = "Metropolitana"
$ws.Cells.Item(104, 4).Value = 44567
$ws.Cells.Item(104, 5).Value = 13
$ws.Cells.Item(104, 6).Value = 100112001
$ws.Cells.Item(104, 7).Value = "Berenjena"
$ws.Cells.Item(104, 8).Value = "Sin especificar"
$ws.Cells.Item(104, 9).Value = "Primera"
$ws.Cells.Item(104, 10).Value = 400
$ws.Cells.Item(104, 11).Value = 9000
$ws.Cells.Item(104, 12).Value = 10000
$ws.Cells.Item(104, 13).Value = 9575
$ws.Cells.Item(104, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(104, 15).Value = "Región Metropolitana"
$ws.Cells.Item(104, 16).Value = 160
$ws.Cells.Item(104, 17).Value = 60
$ws.Cells.Item(104, 18).Value = "Hortaliza"
